# Trade #21 closed at 2026-02-17 08:22:02 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet: refresh aggregate metrics (current capital, P&L, trade
#    counts, win rate) to account for the newly closed trade.
#  - Strategy Status sheet: refresh the MarketMaking strategy row with the
#    same updated figures.
#  - All Trades / MarketMaking sheets: append the new trade record (#21)
#    as row 22.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.26   # Current Capital
$summary.Range("B4").Value = -0.74     # Total P&L $
$summary.Range("B5").Value = -0.7      # Total P&L %
$summary.Range("B6").Value = 21        # Total Trades
$summary.Range("B8").Value = 12        # Losing Trades
$summary.Range("B9").Value = 23.81     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.26      # Capital
$status.Range("D4").Value = 21         # Trades
$status.Range("E4").Value = -0.74      # P&L $
$status.Range("F4").Value = -0.74      # P&L %
$status.Range("G4").Value = 23.81      # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#21 / spreadsheet row 22) to a
# trades-log style worksheet (shared layout between "All Trades" and
# "MarketMaking" sheets).
# ---------------------------------------------------------------------
function Add-TradeRow22($ws) {
    $ws.Range("A22").Value = 21

    # Date / Time columns contain plain text that looks like a date or a
    # time value. A leading apostrophe forces Excel to store them as text
    # instead of auto-converting them into date/time serial numbers.
    $ws.Range("B22").Value = "'2026-02-17"
    $ws.Range("C22").Value = "'08:21:56"

    $ws.Range("D22").Value = "MarketMaking"
    $ws.Range("E22").Value = "UP"
    $ws.Range("F22").Value = 0.71
    $ws.Range("G22").Value = 0.7
    $ws.Range("H22").Value = "CLOSED"
    $ws.Range("I22").Value = -1.4085
    $ws.Range("J22").Value = -0.01
    $ws.Range("K22").Value = 99.26
    $ws.Range("L22").Value = 0
    $ws.Range("M22").Value = 0
    $ws.Range("N22").Value = 0.6
    $ws.Range("O22").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P22").Value = "early_exit"
    $ws.Range("Q22").Value = 0.13
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow22 $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow22 $marketMaking
